$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data previously on row 10 moves to row 9, and the data previously on
# row 9 moves to row 10 (columns A, B, E, F, G, H, Q, R only - the rest of
# each row's values are identical between the two rows already).

# New row 9 values (previously on row 10)
$ws.Range("A9").Value = 111519524
$ws.Range("B9").Value = 77515
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 633583.7615760232
$ws.Range("R9").Value = 7117850.915647855

# New row 10 values (previously on row 9)
$ws.Range("A10").Value = 111519523
$ws.Range("B10").Value = 89686
$ws.Range("E10").Value = 658
$ws.Range("F10").Value = "Rosenticka"
$ws.Range("G10").Value = "Rhodofomes roseus"
$ws.Range("H10").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q10").Value = 633714.5983269843
$ws.Range("R10").Value = 7117626.805168894
